$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.711.62'
$ws.Range("E2").Value = '  +3.94%  '

$ws.Range("D3").Value = '1.923.42'
$ws.Range("E3").Value = '  +2.35%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("B5").Value = 'XRP'
$ws.Range("C5").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D5").Value = '0.704'
$ws.Range("E5").Value = '  +2.99%  '

$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").Value = '249.79'
$ws.Range("E6").Value = '  +1.49%  '

$ws.Range("E7").Value = '  +0.10%  '

$ws.Range("D8").Value = '44.33'
$ws.Range("E8").Value = '  +1.85%  '

$ws.Range("D9").Value = '58.64'
$ws.Range("E9").Value = '  +9.40%  '

$ws.Range("E10").Value = '  +3.89%  '

$ws.Range("D11").Value = '0.0764'
$ws.Range("E11").Value = '  +3.62%  '

$ws.Range("D12").Value = '0.100'
$ws.Range("E12").Value = '  +2.96%  '

$ws.Range("D13").Value = '14.62'
$ws.Range("E13").Value = '  +8.36%  '

$ws.Range("D14").Value = '0.809'
$ws.Range("E14").Value = '  +5.93%  '

$ws.Range("D15").Value = '2.203.66'
$ws.Range("E15").Value = '  +2.42%  '

$ws.Range("E16").Value = '  +4.53%  '

$ws.Range("D17").Value = '1.925.12'
$ws.Range("E17").Value = '  +2.29%  '

$ws.Range("D18").Value = '36.712.79'
$ws.Range("E18").Value = '  +3.90%  '

$ws.Range("D19").Value = '74.45'
$ws.Range("E19").Value = '  +2.29%  '

$ws.Range("E20").Value = '  +4.95%  '

$ws.Range("D21").Value = '252.01'
$ws.Range("E21").Value = '  +3.30%  '

$ws.Range("E22").Value = '  +3.72%  '

$ws.Range("E23").Value = '  +5.42%  '

$ws.Range("D24").Value = '2.67'
$ws.Range("E24").Value = '  +1.96%  '

$ws.Range("E25").Value = '  +0.04%  '

$ws.Range("D26").Value = '2.20'
$ws.Range("E26").Value = '  +0.33%  '

$ws.Range("D27").Value = '168.34'
$ws.Range("E27").Value = '  +1.54%  '

$ws.Range("D28").Value = '8.81'
$ws.Range("E28").Value = '  +3.53%  '

$ws.Range("D29").Value = '18.87'
$ws.Range("E29").Value = '  +3.34%  '

$ws.Range("E30").Value = '  +2.34%  '

$ws.Range("D31").Value = '4.57'
$ws.Range("E31").Value = '  +6.64%  '

$ws.Range("D32").Value = '0.0620'
$ws.Range("E32").Value = '  +4.70%  '

$ws.Range("D33").Value = '1.96'
$ws.Range("E33").Value = '  -3.04%  '

$ws.Range("E34").Value = '  +5.74%  '

$ws.Range("E35").Value = '  +0.10%  '

$ws.Range("D36").Value = '0.0877'
$ws.Range("E36").Value = '  +20.29%  '

$ws.Range("E37").Value = '  -10.87%  '

$ws.Range("D38").Value = '0.897'
$ws.Range("E38").Value = '  +6.93%  '

$ws.Range("D39").Value = '17.77'
$ws.Range("E39").Value = '  +48.82%  '

$ws.Range("D40").Value = '2.04'
$ws.Range("E40").Value = '  +5.24%  '

$ws.Range("D41").Value = '106.67'
$ws.Range("E41").Value = '  +11.14%  '

$ws.Range("E42").Value = '  +5.19%  '

$ws.Range("D43").Value = '17.43'
$ws.Range("E43").Value = '  -1.65%  '

$ws.Range("E44").Value = '  +3.75%  '

$ws.Range("D45").Value = '1.340.61'
$ws.Range("E45").Value = '  +2.91%  '

$ws.Range("D46").Value = '2.59'
$ws.Range("E46").Value = '  +8.76%  '

$ws.Range("E47").Value = '  +1.54%  '

$ws.Range("E48").Value = '  +2.41%  '

$ws.Range("D49").Value = '2.79'
$ws.Range("E49").Value = '  +2.51%  '

$ws.Range("D50").Value = '6.43'
$ws.Range("E50").Value = '  +3.56%  '

$ws.Range("D51").Value = '43.39'
$ws.Range("E51").Value = '  +3.36%  '
